# Fruta / hortaliza, semanal
# This edit reorders the weekly price-report data rows (rows 2-14) on the
# single worksheet. No new data is introduced - the existing 13 rows are
# simply rearranged into a different order. Columns A, B, C, E, F, G, H, I,
# J, K are identical across all rows, so only columns D and L:T actually
# change value as rows move; we rewrite all of them per row for clarity
# and correctness.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D, L, M, N, O, P, Q, R, S, T for rows 2 through 14,
# reflecting the post-edit row order described by the diff.
$rows = @(
    @{ Row = 2;  D = 44917; L = "Segunda"; M = 250; N = 20000; O = 23000; P = 21800; Q = "`$/caja 18 kilos";          R = "Región de Coquimbo";                      S = 1211; T = 18 },
    @{ Row = 3;  D = 44545; L = "Primera"; M = 200; N = 24000; O = 25000; P = 24500; Q = "`$/bandeja 18 kilos";       R = "Región de Coquimbo";                      S = 1361; T = 18 },
    @{ Row = 4;  D = 44881; L = "Segunda"; M = 300; N = 41000; O = 42000; P = 41500; Q = "`$/bandeja 18 kilos";       R = "Región de Coquimbo";                      S = 2306; T = 18 },
    @{ Row = 5;  D = 44524; L = "Segunda"; M = 200; N = 27000; O = 28000; P = 27500; Q = "`$/bandeja 18 kilos";       R = "Provincia de San Felipe de Aconcagua";    S = 1528; T = 18 },
    @{ Row = 6;  D = 44169; L = "Primera"; M = 250; N = 20000; O = 22000; P = 21000; Q = "`$/bandeja 18 kilos";       R = "Provincia de San Felipe de Aconcagua";    S = 1167; T = 18 },
    @{ Row = 7;  D = 44533; L = "Primera"; M = 140; N = 14000; O = 15000; P = 14500; Q = "`$/caja 10 kilos";          R = "Región de O'Higgins";                     S = 1450; T = 10 },
    @{ Row = 8;  D = 44901; L = "Segunda"; M = 200; N = 17000; O = 18000; P = 17500; Q = "`$/bandeja 18 kilos";       R = "Región de O'Higgins";                     S = 972;  T = 18 },
    @{ Row = 9;  D = 44880; L = "Primera"; M = 200; N = 33000; O = 34000; P = 33500; Q = "`$/caja 10 kilos";          R = "Región de O'Higgins";                     S = 3350; T = 10 },
    @{ Row = 10; D = 44174; L = "Primera"; M = 300; N = 19000; O = 20000; P = 19500; Q = "`$/bandeja 18 kilos";       R = "Región Metropolitana";                    S = 1083; T = 18 },
    @{ Row = 11; D = 44895; L = "Segunda"; M = 130; N = 19000; O = 20000; P = 19462; Q = "`$/caja 16 kilos granel";   R = "Región de O'Higgins";                     S = 1216; T = 16 },
    @{ Row = 12; D = 44160; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; Q = "`$/bandeja 18 kilos";       R = "Provincia de San Felipe de Aconcagua";    S = 1361; T = 18 },
    @{ Row = 13; D = 44894; L = "Segunda"; M = 130; N = 19000; O = 20000; P = 19462; Q = "`$/caja 16 kilos granel";   R = "Región de O'Higgins";                     S = 1216; T = 16 },
    @{ Row = 14; D = 44544; L = "Segunda"; M = 250; N = 20000; O = 22000; P = 21000; Q = "`$/bandeja 18 kilos";       R = "Provincia de San Felipe de Aconcagua";    S = 1167; T = 18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $r.R   # R: Origen
    $ws.Cells.Item($row, 19).Value = $r.S   # S: Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $r.T   # T: Kg / unidad
}
